$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (13-21): alexnet / googlenet / resnet152 benchmarks ---
# Row 13 - alexnet (76.8Mhz)
$ws.Range("A13").Value = "alexnet (76.8Mhz)"
$ws.Range("B13").Value = 1.6493560477
$ws.Range("C13").Value = 0.52078262350899995
$ws.Range("D13").Value = 34.441567291299997
$ws.Range("E13").Value = 30.6971039182

# Row 14 - alexnet (537.6Mhz)
$ws.Range("A14").Value = "alexnet (537.6Mhz)"
$ws.Range("B14").Value = 1.88985
$ws.Range("C14").Value = 2.14049821429
$ws.Range("D14").Value = 36.168750000000003
$ws.Range("E14").Value = 32.557142857099997

# Row 15 - alexnet (998.4Mhz)
$ws.Range("A15").Value = "alexnet (998.4Mhz)"
$ws.Range("B15").Value = 2.0312729411800001
$ws.Range("C15").Value = 4.5826635294100004
$ws.Range("D15").Value = 38.451764705899997
$ws.Range("E15").Value = 35.2623529412

# Row 16 - googlenet (76.8Mhz)
$ws.Range("A16").Value = "googlenet (76.8Mhz)"
$ws.Range("B16").Value = 1.67863615561
$ws.Range("C16").Value = 0.57695325269700004
$ws.Range("D16").Value = 37.509398496199999
$ws.Range("E16").Value = 33.602729650199997

# Row 17 - googlenet (537.6Mhz)
$ws.Range("A17").Value = "googlenet (537.6Mhz)"
$ws.Range("B17").Value = 1.85158557692
$ws.Range("C17").Value = 2.7311394230800001
$ws.Range("D17").Value = 39.699038461500002
$ws.Range("E17").Value = 35.922596153800001

# Row 18 - googlenet (998.4Mhz)
$ws.Range("A18").Value = "googlenet (998.4Mhz)"
$ws.Range("B18").Value = 1.9103611111100001
$ws.Range("C18").Value = 6.0862236111100003
$ws.Range("D18").Value = 42.210416666699999
$ws.Range("E18").Value = 39.032638888900003

# Row 19 - resnset152 (76.8Mhz)
$ws.Range("A19").Value = "resnset152 (76.8Mhz)"
$ws.Range("B19").Value = 2.2093594470000002
$ws.Range("C19").Value = 0.59300543120500004
$ws.Range("D19").Value = 40.184002633299997
$ws.Range("E19").Value = 36.210418038199997

# Row 20 - resnet152 (537.6Mhz)
$ws.Range("A20").Value = "resnet152 (537.6Mhz)"
$ws.Range("B20").Value = 2.3617817109099999
$ws.Range("C20").Value = 2.6490668633199999
$ws.Range("D20").Value = 42.306293018700003
$ws.Range("E20").Value = 38.3761061947

# Row 21 - resnet152 (998.4Mhz)
$ws.Range("A21").Value = "resnet152 (998.4Mhz)"
$ws.Range("B21").Value = 2.4679863013699999
$ws.Range("C21").Value = 5.78504657534
$ws.Range("D21").Value = 45.248630136999999
$ws.Range("E21").Value = 41.871232876699999

# --- GPU MEM column (H) ---
$ws.Range("H2").Value = "GPU MEM"
$ws.Range("H13").Value = "720 MB"
$ws.Range("H16").Value = "820 MB"
$ws.Range("H19").Value = "2224 MB"

# --- Update selection to match final state ---
$ws.Range("H21").Select()
